# The document has two logos that each appear twice (once in the
# "first page" header/footer, once in the "default" header/footer):
#   - the Pearson logo, an inline picture inside both footers
#     (docPr/cNvPr name "image1.png")
#   - the BTEC logo, an inline picture inside both headers
#     (docPr/cNvPr name "image2.jpg")
#
# The edit swaps those display names:
#   footers:  image1.png -> image2.png
#   headers:  image2.jpg -> image1.jpg
#
# Note: this runtime's Headers(n)/Footers(n) write path targets the
# physical header/footer part "opposite" the one the same index reads
# from (a quirk of this COM host), so the index used below is the one
# empirically verified (by round-tripping through run_com.py and
# inspecting the saved part XML) to land the change on the intended
# part.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineLogo($headerFooter, $newName) {
    # Go through the shape's own (tight, single-character) Range rather
    # than the whole header/footer Range — setting Name directly off the
    # full HeaderFooter.Range throws ("addressed block not found"), but
    # re-deriving a minimal Range from the shape itself resolves cleanly.
    $shape = $headerFooter.Range.InlineShapes.Item(1)
    $shapeRange = $shape.Range
    $shapeRange.InlineShapes.Item(1).Name = $newName
}

# BTEC logo in the "first page" header (physical header1.xml, docPr id=1)
Rename-InlineLogo $sec.Headers.Item(2) "image1.jpg"

# BTEC logo in the "default" header (physical header2.xml, docPr id=3)
Rename-InlineLogo $sec.Headers.Item(1) "image1.jpg"

# Pearson logo in the "first page" footer (physical footer1.xml, docPr id=2)
Rename-InlineLogo $sec.Footers.Item(2) "image2.png"

# Pearson logo in the "default" footer (physical footer2.xml, docPr id=4)
Rename-InlineLogo $sec.Footers.Item(1) "image2.png"
